{"js": "// Move the \"_GoBack\" bookmark from the end of the document to the run of\n// empty paragraphs between \"Respond to audio input with graphics\" and\n// \"Port Audio Implementation:\", collapsing those three empty paragraphs\n// into a single one that now carries the bookmark.\n\nconst body = context.document.body;\n\n// Remove the old \"_GoBack\" bookmark first (it currently sits at the very\n// end of the document, right after \"...lockable data structure.\"). Doing\n// this before inserting the new one avoids any ambiguity between two\n// same-named bookmarks.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the run of (three) empty paragraphs that sits right after the\n// \"Respond to audio input with graphics\" paragraph and right before the\n// \"Port Audio Implementation:\" paragraph.\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === \"Respond to audio input with graphics\") {\n    anchorIndex = i;\n    break;\n  }\n}\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find 'Respond to audio input with graphics' paragraph\");\n}\n\nlet emptyStart = anchorIndex + 1;\nlet emptyEnd = emptyStart;\nwhile (emptyEnd < items.length && items[emptyEnd].text === \"\") {\n  emptyEnd++;\n}\n// [emptyStart, emptyEnd) are the empty paragraphs (expected 3 of them).\n\n// Keep the first empty paragraph, delete the rest so only one remains.\nfor (let i = emptyEnd - 1; i > emptyStart; i--) {\n  items[i].delete();\n}\nawait context.sync();\n\n// Insert the \"_GoBack\" bookmark into the remaining empty paragraph.\nconst keptParagraph = items[emptyStart];\nkeptParagraph.getRange().insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Move the \"_GoBack\" bookmark from the end of the document to the run of\n# empty paragraphs between \"Respond to audio input with graphics\" and\n# \"Port Audio Implementation:\", collapsing those three empty paragraphs\n# into a single one that now carries the bookmark.\n\n$d = $word.ActiveDocument\n\n# Remove the old \"_GoBack\" bookmark first (it currently sits at the very\n# end of the document, right after \"...lockable data structure.\"). Doing\n# this before inserting the new one avoids any ambiguity between two\n# same-named bookmarks.\nif ($d.Bookmarks.Exists('_GoBack')) {\n    $d.Bookmarks.Item('_GoBack').Delete()\n}\n\n# Find the \"Respond to audio input with graphics\" paragraph, then locate\n# the run of empty paragraphs that immediately follows it.\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -eq \"Respond to audio input with graphics`r\") {\n        $anchorIndex = $i\n        break\n    }\n}\nif ($anchorIndex -eq -1) {\n    throw \"Could not find 'Respond to audio input with graphics' paragraph\"\n}\n\n$emptyStart = $anchorIndex + 1\n$emptyEnd = $emptyStart\nwhile ($emptyEnd -le $d.Paragraphs.Count -and $d.Paragraphs.Item($emptyEnd).Range.Text -eq \"`r\") {\n    $emptyEnd++\n}\n# Empty paragraphs occupy indices [$emptyStart, $emptyEnd) (expected 3 of them).\n\n# Keep the first empty paragraph, delete the rest so only one remains.\nfor ($i = $emptyEnd - 1; $i -gt $emptyStart; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n}\n\n# Insert the \"_GoBack\" bookmark into the remaining empty paragraph.\n$kept = $d.Paragraphs.Item($emptyStart)\n$d.Bookmarks.Add('_GoBack', $kept.Range)\n"}
